# Add a "timeToCookMins" (int, not null) field to the "recipe" table in the
# schema worksheet. The "recipe" table's fields currently occupy rows 8-13:
#   8  recipeID         int           yes   not null
#   9  name              varchar(80)        not null
#   10 description       varchar(280)       not null
#   11 difficultyLevel   int                not null
#   12 img_url           varchar(80)        not null
#   13 author            varchar(80)        not null
#
# The new column is inserted right after "difficultyLevel" (row 11) and
# before "img_url" (old row 12), i.e. as the new row 12 — pushing img_url,
# author, and every following table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 12; this shifts rows 12..65 down to 13..66 and
# grows the outline/merged-cell ranges accordingly.
$ws.Rows.Item(12).Insert()

# Inherit the look (borders/shading/alignment) of a sibling field row
# (row 11, "difficultyLevel") instead of the blank style Insert() leaves.
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(12).OutlineLevel = 1

# Fill in the new field's data: Table Field / Type / Constraints?
# (Primary Key? stays blank, matching the other non-key fields.)
$ws.Range("B12").Value = "timeToCookMins"
$ws.Range("C12").Value = "int"
$ws.Range("E12").Value = "not null"

# Match the saved selection/view state from the edit (cursor sitting on the
# new row's Constraints? cell, no stale scroll position).
$ws.Range("E13").Select()
